# Auto update: 2025-11-29 03:48:30
# Apply targeted numeric updates to the 미장_비트코인_분석 worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Riot Platforms, Inc. / RIOT)
$ws.Range("K2").Value = 57.7
$ws.Range("N2").Value = 85.82376350509293

# Row 3 (Coinbase Global, Inc. / COIN)
$ws.Range("K3").Value = 52.9
$ws.Range("N3").Value = 85.82376350509293

# Row 4 (Bitcoin USD / BTC-USD)
$ws.Range("D4").Value = 90850.27
$ws.Range("F4").Value = 4.66
$ws.Range("K4").Value = 50.7
$ws.Range("N4").Value = 85.82376350509293

# Row 5 (MARA Holdings, Inc. / MARA)
$ws.Range("K5").Value = 50.1
$ws.Range("N5").Value = 85.82376350509293

# Row 6 (Strategy Inc / MSTR)
$ws.Range("K6").Value = 47.7
$ws.Range("N6").Value = 85.82376350509293
